$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Switch the three tables (slides 14, 15, 16) from the custom "Table_0"
#    style to the built-in table style used after the edit.
# ---------------------------------------------------------------------------
$newTableStyleId = "{BD4BAFAE-6608-4B97-BCA3-1017173226D9}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shpIdx = 1; $shpIdx -le $slide.Shapes.Count; $shpIdx++) {
        $shape = $slide.Shapes.Item($shpIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the deck's theme from the "Integral" (Red Violet) palette to
#    the standard "Office Theme" palette.
# ---------------------------------------------------------------------------
function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches MsoThemeColorSchemeIndex 1-12:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = HexToRgb $officeTheme[$i - 1]
}
